$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.295.57'
$ws.Range('E2').Value = '  +2.96%  '
$ws.Range('D3').Value = '2.323.37'
$ws.Range('E3').Value = '  +1.09%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '545.41'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.49%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '130.94'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.00%  '
$ws.Range('E7').Value = '  +0.09%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.579'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -1.17%  '
$ws.Range('D9').Value = '2.321.25'
$ws.Range('E9').Value = '  +1.05%  '
$ws.Range('E10').Value = '  +0.30%  '
$ws.Range('E11').Value = '  -0.02%  '
$ws.Range('E12').Value = '  +0.14%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.335'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.13%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '23.66'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.92%  '
$ws.Range('D15').Value = '60.267.86'
$ws.Range('E15').Value = '  +3.07%  '
$ws.Range('D16').Value = '2.735.52'
$ws.Range('E16').Value = '  +1.03%  '
$ws.Range('E17').Value = '  +0.02%  '
$ws.Range('D18').Value = '2.354.30'
$ws.Range('E18').Value = '  +1.76%  '
$ws.Range('E19').Value = '  -0.27%  '
$ws.Range('E20').Value = '  -1.65%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '313.57'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.81%  '
$ws.Range('E22').Value = '  +0.46%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.998'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.28%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '63.65'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.62%  '
$ws.Range('E25').Value = '  +1.79%  '
$ws.Range('E26').Value = '  +0.05%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.85'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.96%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.35'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +3.47%  '
$ws.Range('B29').Value = 'SuiNetwork'
$ws.Range('C29').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.20'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +8.77%  '
$ws.Range('B30').Value = 'Monero'
$ws.Range('C30').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '173.51'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.31%  '
$ws.Range('E31').Value = '  +1.17%  '
$ws.Range('D32').Value = '0.0₃0731'
$ws.Range('E32').Value = '  +0.36%  '
$ws.Range('E33').Value = '  +1.36%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.37'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +9.83%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.381'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.35%  '
$ws.Range('E36').Value = '  +0.02%  '
$ws.Range('E37').Value = '  -0.49%  '
$ws.Range('E38').Value = '  -0.09%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '4.06'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +1.81%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '323.75'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +10.45%  '
$ws.Range('B41').Value = 'Stacks'
$ws.Range('C41').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.53'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +1.89%  '
$ws.Range('B42').Value = 'OKB'
$ws.Range('C42').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '37.93'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -1.11%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '137.73'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -2.33%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.49'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.75%  '
$ws.Range('E45').Value = '  -1.26%  '
$ws.Range('E46').Value = '  +4.60%  '
$ws.Range('B47').Value = 'Mantle'
$ws.Range('C47').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.560'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.62%  '
$ws.Range('B48').Value = 'Hedera'
$ws.Range('C48').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0494'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.80%  '
$ws.Range('E49').Value = '  +0.80%  '
$ws.Range('E50').Value = '  +13.92%  '
$ws.Range('E51').Value = '  +0.66%  '
